# fix compile for XLUtils
# - Add "Home_Page" values to Test Steps!D2 and Test Steps!D3
# - Move the active cell/selection on "Test Cases" to D2 (no longer the active sheet)
# - Move the active cell/selection on "Test Steps" to E3 (becomes the active sheet)

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestSteps = $wb.Worksheets.Item("Test Steps")

# Populate the new Page Name values on the "Test Steps" sheet
$wsTestSteps.Range("D2").Value = "Home_Page"
$wsTestSteps.Range("D3").Value = "Home_Page"

# Update selection on "Test Cases" (it stops being the active/tabSelected sheet)
$wsTestCases.Select()
$wsTestCases.Range("D2").Select()

# Make "Test Steps" the active sheet and update its selection
$wsTestSteps.Select()
$wsTestSteps.Range("E3").Select()
